# Auto-generated edit script.
#
# The document body contains one 20-row x 5-column table of simple
# addition/subtraction problems (e.g. "21+44=65"). The commit replaces
# every one of the 100 problem strings with a new problem string, in
# table (row-major) order.
#
# NOTE: several old values repeat (e.g. "7+91=98" appears twice, in
# different cells, mapping to two different new values), so a single
# document-wide Find/Replace cannot be used unambiguously. Instead we
# walk the table cell-by-cell (matching the same left-to-right,
# top-to-bottom order the XML diff uses) and overwrite each cell's text
# range directly. We re-resolve each cell's Range on every iteration so
# that earlier replacements (which can change the length of the text)
# never throw off the position of subsequent cells.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Expected existing value per cell, used only as a sanity check.
$oldValues = @(
    "21+44=65",
    "15+50=65",
    "70-69=1",
    "55-28=27",
    "95-69=26",
    "70-33=37",
    "96-22=74",
    "25-0=25",
    "17+40=57",
    "5+21=26",
    "51+4=55",
    "98-67=31",
    "29-17=12",
    "41-11=30",
    "4+24=28",
    "7+91=98",
    "72+5=77",
    "22-12=10",
    "63-29=34",
    "35-12=23",
    "98+0=98",
    "29+28=57",
    "1+5=6",
    "48-43=5",
    "7+91=98",
    "33+34=67",
    "27-23=4",
    "5+1=6",
    "44+34=78",
    "81-57=24",
    "35-17=18",
    "73+2=75",
    "49-35=14",
    "31+28=59",
    "36+59=95",
    "78-56=22",
    "27+59=86",
    "7+67=74",
    "20+7=27",
    "7+76=83",
    "17+76=93",
    "1+90=91",
    "81-2=79",
    "12-10=2",
    "70-55=15",
    "61-21=40",
    "91-39=52",
    "81-74=7",
    "46+8=54",
    "45-43=2",
    "86-5=81",
    "27-21=6",
    "79-6=73",
    "98-68=30",
    "90-41=49",
    "0+89=89",
    "76+20=96",
    "39+42=81",
    "67+0=67",
    "57+19=76",
    "87+1=88",
    "81+3=84",
    "55-50=5",
    "32-30=2",
    "21-18=3",
    "57-50=7",
    "42-20=22",
    "75+13=88",
    "78-30=48",
    "55-34=21",
    "3+6=9",
    "27+13=40",
    "80+13=93",
    "65-14=51",
    "9+88=97",
    "18+12=30",
    "71-11=60",
    "0+36=36",
    "24+19=43",
    "66+29=95",
    "25-23=2",
    "99-52=47",
    "31-3=28",
    "62-22=40",
    "87-32=55",
    "19-8=11",
    "18-4=14",
    "34+40=74",
    "93-63=30",
    "23+63=86",
    "85-70=15",
    "86-22=64",
    "42-3=39",
    "1+43=44",
    "81-32=49",
    "48+29=77",
    "72-16=56",
    "7+70=77",
    "73-63=10",
    "84-15=69"
)

# Replacement value per cell, in the same order.
$newValues = @(
    "12+8=20",
    "5+85=90",
    "61-3=58",
    "40+32=72",
    "79-38=41",
    "34-0=34",
    "22+65=87",
    "61-18=43",
    "46-17=29",
    "30+50=80",
    "9+83=92",
    "0+95=95",
    "23+55=78",
    "76-62=14",
    "90-19=71",
    "97+2=99",
    "41+0=41",
    "17-0=17",
    "74+18=92",
    "4+70=74",
    "90-54=36",
    "90-9=81",
    "96-32=64",
    "20+49=69",
    "45+22=67",
    "48+21=69",
    "99-64=35",
    "19+69=88",
    "31+68=99",
    "40+23=63",
    "9+6=15",
    "93-40=53",
    "92+7=99",
    "98-81=17",
    "25+52=77",
    "23+50=73",
    "71+24=95",
    "17+3=20",
    "58-32=26",
    "94-62=32",
    "72+10=82",
    "29+26=55",
    "30+1=31",
    "87-30=57",
    "64-20=44",
    "14-4=10",
    "61-36=25",
    "44+25=69",
    "50+26=76",
    "43+20=63",
    "51-48=3",
    "27+40=67",
    "60-28=32",
    "36+60=96",
    "97-14=83",
    "74+0=74",
    "3+45=48",
    "22+13=35",
    "74-56=18",
    "36-0=36",
    "67-26=41",
    "97-4=93",
    "35+32=67",
    "8+31=39",
    "26+13=39",
    "81-3=78",
    "4+2=6",
    "47-26=21",
    "78-20=58",
    "69-28=41",
    "69-4=65",
    "73-60=13",
    "19-18=1",
    "71-5=66",
    "98-37=61",
    "24+12=36",
    "34-21=13",
    "68+14=82",
    "4+74=78",
    "57+30=87",
    "36-18=18",
    "23+70=93",
    "40+22=62",
    "17-6=11",
    "54+22=76",
    "33-6=27",
    "76+19=95",
    "24-21=3",
    "52-39=13",
    "44+51=95",
    "1+94=95",
    "71-35=36",
    "42-41=1",
    "18-14=4",
    "38+0=38",
    "23+25=48",
    "36-23=13",
    "91-74=17",
    "99-11=88",
    "53-53=0"
)

$rows = 20
$cols = 5
$idx = 0
$mismatches = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cellRange = $cell.Range
        # Trim the trailing end-of-cell mark so we only touch the
        # visible text, not the cell's paragraph/cell-mark characters.
        $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)

        $expectedOld = $oldValues[$idx]
        $newVal = $newValues[$idx]

        if ($textRange.Text -ne $expectedOld) {
            $mismatches = $mismatches + 1
        }

        $textRange.Text = $newVal

        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells ($mismatches unexpected)"
